$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 (OTROS group): update VENTA and POR CUMPLIR
$ws.Range("D2").Value = 3059.12
$ws.Range("E2").Value = -3059.12

# Row 4 (TOTAL row): update VENTA, POR CUMPLIR and CUMPLIMIENTO (%)
$ws.Range("D4").Value = 3307.6
$ws.Range("E4").Value = 14192.4
$ws.Range("F4").Value = 0.1890057142857143
